# Update "想去人数" (want-to-go count) values across sheets, reflecting the
# newer snapshot of the data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (F column = want-to-go count) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 9082
$wsExhibit.Range("F3").Value = 101
$wsExhibit.Range("F4").Value = 242
$wsExhibit.Range("F5").Value = 111
$wsExhibit.Range("F6").Value = 1510
$wsExhibit.Range("F7").Value = 1422
$wsExhibit.Range("F8").Value = 252
$wsExhibit.Range("F10").Value = 338
$wsExhibit.Range("F11").Value = 92

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 17

# --- Sheet "全部类型" (combined view, includes 演出 row inserted at row 9) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9082
$wsAll.Range("F3").Value = 101
$wsAll.Range("F4").Value = 242
$wsAll.Range("F5").Value = 111
$wsAll.Range("F6").Value = 1510
$wsAll.Range("F7").Value = 1422
$wsAll.Range("F8").Value = 252
$wsAll.Range("F9").Value = 17
$wsAll.Range("F11").Value = 338
$wsAll.Range("F12").Value = 92
